$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Accipiter gentilis"
$ws.Range("B7").Value = 9.7
$ws.Range("C7").Value = 0

$ws.Range("D12").Select()
